# Updated the exception handling
# - Row 2: candidate became a "Fresher" (no longer "Experienced"), so experience
#   years / notice period are cleared, email/phone updated.
# - Row 3: the second test candidate row is cleared out entirely.
# - The hyperlink on G3 (tied to the removed row) is deleted, and the
#   hyperlink on G2 is repointed at the updated e-mail address.
# - Selection moves to G2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 values ---
$ws.Range("G2").Value = "George.Thompson1a741@gmail.com"
$ws.Range("H2").Value = 1111111111
$ws.Range("I2").Value = "Fresher"
$ws.Range("J2").ClearContents()
$ws.Range("L2").ClearContents()

# --- Clear row 3 entirely (keep formatting) ---
$ws.Range("A3:P3").ClearContents()

# --- Fix up hyperlinks: repoint G2, remove G3's ---
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$G$2') {
        $hl.Address = "mailto:George.Thompson1a741@gmail.com"
    }
}
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$G$3') {
        $hl.Delete()
    }
}

# --- Update selection to G2 ---
[void]$ws.Range("G2").Select()
